$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-11-05 Tuesday" "2024-11-06 Wednesday"

Replace-Text "294×9=" "205×2="
Replace-Text "446×4=" "486×4="
Replace-Text "805×5=" "573×8="
Replace-Text "739×2=" "495×7="
Replace-Text "736×3=" "651×2="

Replace-Text "631×4=" "821×4="
Replace-Text "809×4=" "801×9="
Replace-Text "578×4=" "442×7="
Replace-Text "831×3=" "309×4="
Replace-Text "335×6=" "855×2="

Replace-Text "993×5=" "456×2="
Replace-Text "464×9=" "927×2="
Replace-Text "947×4=" "492×8="
Replace-Text "855×6=" "231×3="
Replace-Text "724×2=" "804×9="

Replace-Text "404×4=" "282×5="
Replace-Text "525×3=" "254×4="
Replace-Text "865×5=" "900×5="
Replace-Text "283×4=" "105×8="
Replace-Text "900×8=" "620×4="

Replace-Text "153×2=" "295×9="
Replace-Text "287×8=" "230×6="
Replace-Text "435×5=" "369×7="
Replace-Text "364×4=" "311×9="
Replace-Text "691×6=" "614×7="
